$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mapping of row -> (DAMSLTag, DialogAct)
$updates = @{
    9   = @("b",  "Acknowledge (Backchannel)")
    15  = @("ba", "Appreciation")
    16  = @("%",  "Uninterpretable")
    22  = @("sd", "Statement-non-opinion")
    28  = @("aa", "Agree/Accept")
    69  = @("sv", "Statement-opinion")
    82  = @("aa", "Agree/Accept")
    102 = @("sv", "Statement-opinion")
    135 = @("aa", "Agree/Accept")
    136 = @("aa", "Agree/Accept")
    138 = @("aa", "Agree/Accept")
    139 = @("ba", "Appreciation")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}

$wb.Save()
